$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the combined header style (bold + thin border + center/top align) ---
# A2 already carries the workbook's sole "bold" style (from the original
# "World" cell), so mutating it in place re-purposes that exact style slot
# into the full combo style instead of allocating a brand new one.
$a2 = $ws.Range("A2")
$a2.Borders.LineStyle = 1
$a2.Borders.Weight = 2
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160

# Copy that now-upgraded style onto A1 (the future header cell).
$a2.Copy()
$ws.Range("A1").PasteSpecial(-4122)

# --- Cell content ---
$ws.Range("A1").Value = "Data"

$a2.ClearFormats()
$a2.Value = 10

$ws.Range("A3").Value = 20
$ws.Range("A4").Value = 30
$ws.Range("A5").Value = 20
$ws.Range("A6").Value = 15
$ws.Range("A7").Value = 30
$ws.Range("A8").Value = 45
